$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure D and E columns stay as text (matches original inlineStr formatting),
# since some numeric-looking strings (e.g. "556.45") would otherwise be
# auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '60.351.96'
$ws.Range('E2').Value = '  +4.22%  '
$ws.Range('D3').Value = '2.434.47'
$ws.Range('E3').Value = '  +3.24%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '556.45'
$ws.Range('E5').Value = '  +2.43%  '
$ws.Range('D6').Value = '139.52'
$ws.Range('E6').Value = '  +3.66%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +3.20%  '
$ws.Range('E9').Value = '  +4.94%  '
$ws.Range('E10').Value = '  +3.96%  '
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('D13').Value = '25.04'
$ws.Range('E13').Value = '  +5.48%  '
$ws.Range('D14').Value = '2.866.08'
$ws.Range('E14').Value = '  +3.22%  '
$ws.Range('D15').Value = '60.291.82'
$ws.Range('E15').Value = '  +4.17%  '
$ws.Range('E16').Value = '  +4.32%  '
$ws.Range('D17').Value = '2.428.77'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('E18').Value = '  +5.84%  '
$ws.Range('D19').Value = '4.43'
$ws.Range('E19').Value = '  +3.07%  '
$ws.Range('D20').Value = '334.59'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '65.36'
$ws.Range('E23').Value = '  +4.31%  '
$ws.Range('D25').Value = '8.63'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = '0.0₃0792'
$ws.Range('E28').Value = '  +6.92%  '
$ws.Range('D29').Value = '1.79'
$ws.Range('E29').Value = '  +2.21%  '
$ws.Range('D30').Value = '6.34'
$ws.Range('E30').Value = '  +3.30%  '
$ws.Range('D31').Value = '169.20'
$ws.Range('E31').Value = '  -0.83%  '
$ws.Range('E32').Value = '  +0.89%  '
$ws.Range('D33').Value = '18.79'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('E35').Value = '  +6.32%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '0.423'
$ws.Range('E39').Value = '  +11.42%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '39.88'
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('D41').Value = '321.64'
$ws.Range('E41').Value = '  +11.40%  '
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('D43').Value = '141.30'
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('E44').Value = '  +3.65%  '
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('D46').Value = '19.64'
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').Value = '0.416'
$ws.Range('E47').Value = '  +7.97%  '
$ws.Range('D48').Value = '0.575'
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('D49').Value = '0.0228'
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('D50').Value = '17.97'
$ws.Range('E50').Value = '  +2.85%  '
$ws.Range('E51').Value = '  -0.20%  '
